$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''27.190.27'
$ws.Range('E2').Value = '  -2.47%  '
$ws.Range('D3').Value = '''1.873.05'
$ws.Range('E3').Value = '  -1.83%  '
$ws.Range('D4').Value = '''1.000'
$ws.Range('E4').Value = '  -0.31%  '
$ws.Range('D5').Value = '''307.45'
$ws.Range('E5').Value = '  -2.09%  '
$ws.Range('D6').Value = '''1.0000'
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('D7').Value = '''0.5141'
$ws.Range('E7').Value = '  +2.47%  '
$ws.Range('E8').Value = '  -1.56%  '
$ws.Range('D9').Value = '''0.07172'
$ws.Range('E9').Value = '  -1.67%  '
$ws.Range('D10').Value = '''0.8889'
$ws.Range('E10').Value = '  -2.30%  '
$ws.Range('D11').Value = '''20.74'
$ws.Range('E11').Value = '  -0.70%  '
$ws.Range('D12').Value = '''0.07587'
$ws.Range('E12').Value = '  -1.07%  '
$ws.Range('D13').Value = '''1.863.27'
$ws.Range('E13').Value = '  -2.62%  '
$ws.Range('E14').Value = '  -2.51%  '
$ws.Range('E15').Value = '  -2.68%  '
$ws.Range('D16').Value = '''1.000'
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('D17').Value = '''0.000008559'
$ws.Range('E17').Value = '  -1.87%  '
$ws.Range('D18').Value = '''14.18'
$ws.Range('E18').Value = '  -2.74%  '
$ws.Range('E19').Value = '  -0.24%  '
$ws.Range('D20').Value = '''27.224.78'
$ws.Range('E20').Value = '  -2.50%  '
$ws.Range('D21').Value = '''5.071'
$ws.Range('E21').Value = '  -1.93%  '
$ws.Range('D22').Value = '''2.084.22'
$ws.Range('E22').Value = '  -4.08%  '
$ws.Range('E23').Value = '  -1.91%  '
$ws.Range('D24').Value = '''6.497'
$ws.Range('E24').Value = '  -1.35%  '
$ws.Range('D25').Value = '''151.03'
$ws.Range('E25').Value = '  -1.98%  '
$ws.Range('D26').Value = '''1.845'
$ws.Range('E26').Value = '  -1.87%  '
$ws.Range('D27').Value = '''18.04'
$ws.Range('E27').Value = '  -1.95%  '
$ws.Range('D28').Value = '''2.139'
$ws.Range('E28').Value = '  -3.81%  '
$ws.Range('D29').Value = '''112.80'
$ws.Range('E29').Value = '  -2.18%  '
$ws.Range('D30').Value = '''4.761'
$ws.Range('E30').Value = '  -3.15%  '
$ws.Range('D31').Value = '''4.700'
$ws.Range('E31').Value = '  +0.97%  '
$ws.Range('D32').Value = '''0.08997'
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range('D33').Value = '''0.05162'
$ws.Range('E33').Value = '  -1.71%  '
$ws.Range('D34').Value = '''3.103'
$ws.Range('E34').Value = '  -3.46%  '
$ws.Range('D35').Value = '''0.7553'
$ws.Range('E35').Value = '  -1.05%  '
$ws.Range('D36').Value = '''1.176'
$ws.Range('E36').Value = '  -4.42%  '
$ws.Range('D37').Value = '''0.02046'
$ws.Range('E37').Value = '  -0.80%  '
$ws.Range('D38').Value = '''2.536'
$ws.Range('E38').Value = '  -0.27%  '
$ws.Range('D39').Value = '''3.034'
$ws.Range('E39').Value = '  +0.37%  '
$ws.Range('D40').Value = '''1.079'
$ws.Range('E40').Value = '  -1.21%  '
$ws.Range('D41').Value = '''0.5373'
$ws.Range('E41').Value = '  -3.73%  '
$ws.Range('D42').Value = '''6.658'
$ws.Range('E42').Value = '  -3.95%  '
$ws.Range('E43').Value = '  +3.18%  '
$ws.Range('D44').Value = '''8.578'
$ws.Range('E44').Value = '  +0.94%  '
$ws.Range('D45').Value = '''0.1486'
$ws.Range('E45').Value = '  -1.72%  '
$ws.Range('D46').Value = '''0.4686'
$ws.Range('E46').Value = '  -3.17%  '
$ws.Range('D47').Value = '''0.9998'
$ws.Range('E47').Value = '  -0.36%  '
$ws.Range('E48').Value = '  -3.94%  '
$ws.Range('D49').Value = '''1.576'
$ws.Range('E49').Value = '  -3.26%  '
$ws.Range('D50').Value = '''65.20'
$ws.Range('E50').Value = '  -3.48%  '
$ws.Range('D51').Value = '''36.56'
$ws.Range('E51').Value = '  -1.26%  '
